# convert time second to minute
#
# D2 (counterOut) and G2:G13 (runtime) were stored as plain numeric-looking
# text in the shared-string table. Re-stamp them as text (via a temporary
# "@" / text number format) so Excel doesn't silently coerce the new
# values back into numbers, then drop the format back to Normal so no
# stray cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetCells = @("D2", "G2", "G3", "G4", "G5", "G6", "G7", "G8", "G9", "G10", "G11", "G12", "G13")

foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "548"

$ws.Range("G2").Value = "8.0"
$ws.Range("G3").Value = "0.0"
$ws.Range("G4").Value = "0.0"
$ws.Range("G5").Value = "0.0"
$ws.Range("G6").Value = "0.0"
$ws.Range("G7").Value = "0.0"
$ws.Range("G8").Value = "0.0"
$ws.Range("G9").Value = "0.0"
$ws.Range("G10").Value = "0.0"
$ws.Range("G11").Value = "0.0"
$ws.Range("G12").Value = "0.0"
$ws.Range("G13").Value = "0.0"

foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
